$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,5,8,9,12,14,16,21,22,25,29,30,32,34,37,39,41,42,43,47,48,50,54,56,60,64,67,72,76,79,83,85,90,91,96,102,103,106,111,117,119,125,127,128,135,137,140,146,153,154,160,161,162,167,173,174,175,176,178,179,184,186,190,193)
foreach ($r in $rows) {
    $ws.Range("G" + $r + ":N" + $r).Value = "N/A"
}
